$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 73

# Columns A, B, D contain values that look like dates/times/numbers but must
# stay as plain text, matching the style of the existing rows above. Prefix
# with an apostrophe to force text entry, then clear the resulting
# "quote prefix" cell format so the cells end up unstyled (default format),
# just like the rest of the data rows.
$ws.Cells.Item($row, 1).Value = "'2023-06-24"
$ws.Cells.Item($row, 2).Value = "'20:39:39"
$ws.Cells.Item($row, 3).Value = "Saturday"
$ws.Cells.Item($row, 4).Value = "'25"

$ws.Range("A73:D73").ClearFormats()

$ws.Cells.Item($row, 5).Value = 122611
$ws.Cells.Item($row, 6).Value = 134078
$ws.Cells.Item($row, 7).Value = 163063
$ws.Cells.Item($row, 8).Value = 133471
$ws.Cells.Item($row, 9).Value = 177620
$ws.Cells.Item($row, 10).Value = 115838
$ws.Cells.Item($row, 11).Value = 203141
$ws.Cells.Item($row, 12).Value = 226141
$ws.Cells.Item($row, 13).Value = 175540
$ws.Cells.Item($row, 14).Value = 104218
$ws.Cells.Item($row, 15).Value = 39568
$ws.Cells.Item($row, 16).Value = 33818
$ws.Cells.Item($row, 17).Value = 52001
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36135
$ws.Cells.Item($row, 20).Value = -1
